$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column F ("Resolution") entirely, shifting G:M left to F:L
$ws.Range("F1:F2").EntireColumn.Delete() | Out-Null

# Update header row
$ws.Range("E1").Value = "Split"

# Update data row values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 5.795400282503119
$ws.Range("C2").Value = "lst, no label"
$ws.Range("E2").Value = "random"
$ws.Range("F2").Value = 955
$ws.Range("G2").Value = 32
$ws.Range("H2").Value = 100
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "2024-11-13"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "11:32:36"
$ws.Range("K2").Value = 466.92
